# Insert a new weekly price record for "Vega Modelo de Temuco - Tuna" at
# row 64, shifting the existing rows 64..76 down to 65..77.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 64 (and everything below it) down by one row.
$ws.Rows.Item(64).Insert()

# Populate the newly-inserted row 64 with the new weekly observation.
$ws.Cells.Item(64, 1).Value  = 10
$ws.Cells.Item(64, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(64, 3).Value  = "La Araucanía"
$ws.Cells.Item(64, 4).Value  = 44995
$ws.Cells.Item(64, 5).Value  = 9
$ws.Cells.Item(64, 6).Value  = "Fruta"
$ws.Cells.Item(64, 7).Value  = 100107
$ws.Cells.Item(64, 8).Value  = "Otros"
$ws.Cells.Item(64, 9).Value  = 100107011
$ws.Cells.Item(64, 10).Value = "Tuna"
$ws.Cells.Item(64, 11).Value = "Sin especificar"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 80
$ws.Cells.Item(64, 14).Value = 20000
$ws.Cells.Item(64, 15).Value = 20000
$ws.Cells.Item(64, 16).Value = 20000
$ws.Cells.Item(64, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(64, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(64, 19).Value = 1250
$ws.Cells.Item(64, 20).Value = 16
